$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell reference -> new text value (Price and Volume(1h) columns)
$updates = @{
    'D2' = '69.720.00'
    'E2' = '  -0.09%  '
    'D3' = '3.679.39'
    'E3' = '  -0.62%  '
    'D4' = '1.00'
    'E4' = '  -0.02%  '
    'D5' = '651.36'
    'E5' = '  -3.95%  '
    'D6' = '161.01'
    'E6' = '  -0.84%  '
    'E7' = '  -0.02%  '
    'D8' = '0.497'
    'E8' = '  -0.04%  '
    'E10' = '  +0.45%  '
    'D11' = '0.441'
    'E11' = '  -0.56%  '
    'E12' = '  -2.22%  '
    'D13' = '4.302.28'
    'E13' = '  -0.66%  '
    'D14' = '32.66'
    'E14' = '  -0.70%  '
    'D15' = '3.664.07'
    'E15' = '  -1.05%  '
    'D16' = '69.743.25'
    'E16' = '  -0.02%  '
    'E17' = '  +0.75%  '
    'D18' = '6.52'
    'E18' = '  +0.12%  '
    'D19' = '15.91'
    'E19' = '  -1.43%  '
    'D20' = '10.36'
    'E20' = '  +5.38%  '
    'D21' = '470.58'
    'E21' = '  -0.46%  '
    'D22' = '0.653'
    'E22' = '  -0.12%  '
    'D23' = '79.75'
    'E23' = '  -0.98%  '
    'D24' = '3.826.87'
    'E24' = '  -0.66%  '
    'D25' = '0.999'
    'E26' = '  -1.53%  '
    'D27' = '11.12'
    'E27' = '  +0.46%  '
    'E28' = '  -4.23%  '
    'E29' = '  -2.43%  '
    'E30' = '  -3.36%  '
    'E31' = '  -2.06%  '
    'E32' = '  -0.02%  '
    'E33' = '  +0.66%  '
    'D34' = '26.72'
    'E34' = '  -0.68%  '
    'D35' = '6.41'
    'E35' = '  -3.22%  '
    'D36' = '3.675.48'
    'E36' = '  -0.43%  '
    'D37' = '8.37'
    'E37' = '  -2.34%  '
    'D39' = '5.87'
    'E39' = '  -5.19%  '
    'D40' = '178.12'
    'E40' = '  +5.22%  '
    'D41' = '1.00'
    'E41' = '  -0.03%  '
    'D42' = '0.0893'
    'E42' = '  -1.75%  '
    'D43' = '2.17'
    'E43' = '  -2.43%  '
    'D44' = '0.929'
    'E44' = '  -1.75%  '
    'D45' = '46.80'
    'E45' = '  -0.45%  '
    'D46' = '29.12'
    'E46' = '  +3.77%  '
    'E47' = '  -1.02%  '
    'E48' = '  -4.77%  '
    'E49' = '  -0.89%  '
    'E50' = '  -4.77%  '
    'E51' = '  -5.96%  '
}

foreach ($cellRef in $updates.Keys) {
    $cell = $ws.Range($cellRef)
    # Force text format so numeric-looking strings (e.g. "1.00", "0.497")
    # keep their exact textual representation instead of becoming numbers.
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$cellRef]
}
